$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("issues")

$newText = "* Load some workspaces, try 5 from different instruments`n* Show detectors on these workspaces`n* Check Instrument Viewer for these workspaces, follow  documentation for  [Instrument Viewer Widget](https://docs.mantidproject.org/workbench/instrumentviewer.html) and some more detail on the picking tools in the [Basic Course](https://docs.mantidproject.org/tutorials/mantid_basic_course/connecting_data_to_instruments/03_investigating_data.html)`n  - [ ] Render tab `n  - [ ] Draw tab `n  - [ ] Pick tab `n  - [ ] Instrument tree tab `n* Test 5 - 10 random algorithms, follow [how to run Algorithms](http://docs.mantidproject.org/nightly/tutorials/mantid_basic_course/algorithms_workspaces_and_histories/02_algorithms.html#algorithms) and [pick your favourite](https://docs.mantidproject.org/algorithms/)`n - [ ] Help button `n - [ ] Validation of inputs `n - [ ] Overwriting workspace `n - [ ] Run the algorithm `n- [ ] Check right-click options of the workspaces, they are listed [in the image in the documentation](https://docs.mantidproject.org/nightly/workbench/workspacetreewidget.html#workbenchworkspacetoolbox).`n"

$ws.Range("C7").Value = $newText
$ws.Rows.Item(7).RowHeight = 99.75
